{"js": "// Applies the \"various help text and comments updated\" edit:\n//  1) Inserts a new \"Default Settings at Plugin Activation Time\" section\n//     (Heading1 + 4 body paragraphs + a blank paragraph) right before the\n//     existing \"PayPal API Code Related\" Heading1 paragraph.\n//  2) Refreshes the \"When the transaction completes on the popup...\"\n//     paragraph so the stale <w:lastRenderedPageBreak/> marker left over\n//     from the old layout is dropped (the visible text is unchanged).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"PayPal API Code Related\" heading paragraph - the new section\n// is inserted immediately before it.\nlet headingPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"PayPal API Code Related\") {\n    headingPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (headingPara) {\n  // Insert in reverse order, always anchoring \"before\" the heading so the\n  // final reading order matches the diff. New paragraphs inherit the\n  // anchor's style (Heading1), so explicitly reset the body ones back to\n  // Normal.\n  const blank = headingPara.insertParagraph(\"\", Word.InsertLocation.before);\n  blank.styleBuiltIn = Word.Style.normal;\n\n  const getDefaultsPara = blank.insertParagraph(\"Main::get_defaults()\", Word.InsertLocation.before);\n  getDefaultsPara.styleBuiltIn = Word.Style.normal;\n\n  const introPara2 = getDefaultsPara.insertParagraph(\n    \"The following function can be used to set the default settings values:\",\n    Word.InsertLocation.before\n  );\n  introPara2.styleBuiltIn = Word.Style.normal;\n\n  const singleActivatePara = introPara2.insertParagraph(\"Main::single_activate()\", Word.InsertLocation.before);\n  singleActivatePara.styleBuiltIn = Word.Style.normal;\n\n  const introPara1 = singleActivatePara.insertParagraph(\n    \"#) The following function is run at plugin activation\",\n    Word.InsertLocation.before\n  );\n  introPara1.styleBuiltIn = Word.Style.normal;\n\n  const titlePara = introPara1.insertParagraph(\"Default Settings at Plugin Activation Time\", Word.InsertLocation.before);\n  titlePara.styleBuiltIn = Word.Style.heading1;\n\n  await context.sync();\n}\n\n// Remove the stale lastRenderedPageBreak marker on the \"When the\n// transaction completes...\" paragraph by re-typing its (unchanged) text.\nconst searchResults = body.search(\"When the transaction completes on the popup\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const hitRange = searchResults.items[0];\n  const paragraph = hitRange.paragraphs.getFirst();\n  paragraph.load(\"text\");\n  await context.sync();\n\n  const fullText = paragraph.text;\n  paragraph.clear();\n  paragraph.insertText(fullText, Word.InsertLocation.start);\n  await context.sync();\n}\n", "ps1": "# Applies the \"various help text and comments updated\" edit:\n#  1) Inserts a new \"Default Settings at Plugin Activation Time\" section\n#     (Heading1 + 4 body paragraphs + a blank paragraph) right before the\n#     existing \"PayPal API Code Related\" Heading1 paragraph.\n#  2) Refreshes the \"When the transaction completes on the popup...\"\n#     paragraph so the stale lastRenderedPageBreak marker left over from\n#     the old layout is dropped (the visible text is unchanged).\n\n$d = $word.ActiveDocument\n\n# --- 1) Insert the new \"Default Settings at Plugin Activation Time\" section ---\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $txt = $d.Paragraphs($i).Range.Text\n  $txt = $txt -replace \"[`r`n\\x07]\", \"\"\n  if ($txt -eq \"PayPal API Code Related\") {\n    $targetIndex = $i\n    break\n  }\n}\n\nif ($targetIndex -ne -1) {\n  $targetRange = $d.Paragraphs($targetIndex).Range\n  for ($n = 0; $n -lt 6; $n++) {\n    $targetRange.InsertParagraphBefore()\n  }\n\n  $lines = @(\n    \"Default Settings at Plugin Activation Time\",\n    \"#) The following function is run at plugin activation\",\n    \"Main::single_activate()\",\n    \"The following function can be used to set the default settings values:\",\n    \"Main::get_defaults()\",\n    \"\"\n  )\n  $styles = @(\"Heading 1\", \"Normal\", \"Normal\", \"Normal\", \"Normal\", \"Normal\")\n\n  for ($n = 0; $n -lt 6; $n++) {\n    $p = $d.Paragraphs($targetIndex + $n)\n    $p.Range.Text = $lines[$n]\n    $p.Range.Style = $d.Styles($styles[$n])\n  }\n}\n\n# --- 2) Re-type the \"When the transaction completes...\" paragraph so the\n#        stale lastRenderedPageBreak marker is dropped ---\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"When the transaction completes on the popup: JS code\"\n$found = $find.Execute()\n\nif ($found) {\n  $para = $rng.Paragraphs(1)\n  $paraRange = $para.Range\n  # Exclude the trailing paragraph mark so we don't merge with the next paragraph.\n  $paraRange.MoveEnd(1, -1) | Out-Null\n  $originalText = $paraRange.Text\n  $paraRange.Text = \"\"\n  $paraRange.InsertAfter($originalText)\n}\n"}
